# Starting to add the earth heroes on hero selected panel.
# Remove the "King Guava" / "Queen Bittergourd" / "Insect Queen (Angela)"
# placeholder rows (name + description), change the Blaster Melon
# description, and update the active sheet selection/view + un-minimize
# the workbook window.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Clear the "King Guava" character row (name + description) but keep formatting.
$ws.Range("A3:B3").ClearContents()

# Clear the "Queen Bittergourd" character row (name + description) but keep formatting.
$ws.Range("A4:B4").ClearContents()

# Clear the "Insect Queen (Angela)" character row (name + description) but keep formatting.
$ws.Range("A27:B27").ClearContents()

# Update the Blaster Melon description text.
$ws.Range("B18").Value = "This dude  got some big seeds to talk about !"

# Scroll the view and move the selection.
$excel.Application.Goto($ws.Range("B1"), $false) | Out-Null
$ws.Range("B18").Select() | Out-Null

# Restore the window from minimized.
$excel.Windows.Item(1).WindowState = -4143 | Out-Null
